# Updated simulation files with Holden scheme
# ---------------------------------------------------------------------------
# The sheet holds a 21-scheme x 18-HKL grid. Four new "Holden" sampling
# schemes (Holden2.5 / Holden5 / Holden10 / Holden15) are inserted into the
# scheme list right after the "Offset*" group and before "HexGrid-*", and
# the duplicated U:AD block (columns 19-28, a stray repeat of the last
# 10 HKL headers) is dropped. The net effect: rows 3-19 keep their original
# row labels, the header row / HKL-label row get fixed up, the leftover
# duplicate columns disappear, and four brand-new rows are appended at the
# bottom of the table for the schemes that used to occupy rows 16-19 before
# the Holden rows were spliced in ahead of them.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the stray duplicate block U1:AD2 (columns 19-28). Clear() removes
#    the cells outright (not just their content) so the used range shrinks
#    back down to column T, matching dimension A1:T23 once rows 20-23 exist.
$ws.Range("U1:AD2").Clear()

# 2) Full logical contents of A1:T23 after the edit (row 1 = column index
#    header, row 2 = HKL-label header, rows 3-23 = one scheme per row).
#    $null entries are left untouched (A1/A... has no cell at all).
$grid = @(
  @($null, 0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18),
  @(0, "HKL", "[4, 0, 0]", "[2, 1, 1]", "[2, 2, 0]", "[2, 0, 0]", "[2, 2, 2]", "[3, 1, 0]", "[1, 1, 0]", "[3, 2, 1]", "1Pair-A", "1Pair-B", "2Pairs-A", "2Pairs-B", "3Pairs-A", "3Pairs-B", "3Pairs-C", "4Pairs", "5A4F", "MaxUnique"),
  @(1, "BT8Hex_2.5", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(2, "BT8Hex_5", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(3, "BT8Hex_10", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(4, "BT8Hex_15", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(5, "Spiral2.5", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(6, "Spiral5", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(7, "Spiral7.5", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(8, "Spiral10", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(9, "Spiral15", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(10, "OffsetF45", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(11, "OffsetA45", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(12, "OffsetFTD", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(13, "OffsetATD", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(14, "Holden2.5", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(15, "Holden5", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(16, "Holden10", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(17, "Holden15", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(18, "HexGrid-90degTilt2.5degRes", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(19, "HexGrid-90degTilt5degRes", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(20, "HexGrid-90degTilt10degRes", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(21, "HexGrid-90degTilt15degRes", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($r = 0; $r -lt $grid.Length; $r++) {
    $rowNum = $r + 1
    $rowData = $grid[$r]
    for ($ci = 0; $ci -lt $rowData.Length; $ci++) {
        $val = $rowData[$ci]
        if ($null -eq $val) { continue }
        $addr = $cols[$ci] + $rowNum
        $ws.Range($addr).Value = $val
    }
}

# 3) New rows 20-23 need the same bold/border/center-top style ("style 1")
#    already used on column A and row 1. Copy formatting from the existing
#    A-column cells above (A16:A19) so the stylesheet is reused as-is
#    instead of minting a near-duplicate style.
$ws.Range("A16:A19").Copy() | Out-Null
$ws.Range("A20:A23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
